$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 170061253.5450011
$ws.Range("C2").Value = -402048568.8690004
$ws.Range("D2").Value = 53
$ws.Range("E2").Value = "T"
$ws.Range("B3").Value = 215751600.2529267
$ws.Range("C3").Value = -466156325.3815762
$ws.Range("D3").Value = 53
$ws.Range("E3").Value = "T"
$ws.Range("B4").Value = 190157095.1846463
$ws.Range("C4").Value = -415950649.4723316
$ws.Range("D4").Value = 53
$ws.Range("E4").Value = "T"
$ws.Range("B5").Value = 170860992.9712256
$ws.Range("C5").Value = -405155427.0714723
$ws.Range("D5").Value = 53
$ws.Range("E5").Value = "T"
$ws.Range("B6").Value = 170667185.431315
$ws.Range("C6").Value = -405104104.4883425
$ws.Range("D6").Value = 53
$ws.Range("E6").Value = "T"
$ws.Range("B7").Value = 198148318.8949276
$ws.Range("C7").Value = -446478395.0038871
$ws.Range("D7").Value = 53
$ws.Range("E7").Value = "T"
$ws.Range("B8").Value = 170922082.497805
$ws.Range("C8").Value = -404787896.7358894
$ws.Range("D8").Value = 53
$ws.Range("E8").Value = "T"
$ws.Range("B9").Value = 190811154.076052
$ws.Range("C9").Value = -418967515.8208964
$ws.Range("D9").Value = 53
$ws.Range("E9").Value = "T"
$ws.Range("B10").Value = 167063073.7771204
$ws.Range("C10").Value = -399736333.0616241
$ws.Range("D10").Value = 53
$ws.Range("E10").Value = "T"
$ws.Range("B11").Value = 191182675.188347
$ws.Range("C11").Value = -417313437.041097
$ws.Range("D11").Value = 53
$ws.Range("E11").Value = "T"
$ws.Range("B12").Value = 226125332.4069945
$ws.Range("C12").Value = -465476449.0371869
$ws.Range("D12").Value = 53
$ws.Range("E12").Value = "T"
$ws.Range("B13").Value = 207534924.9715279
$ws.Range("C13").Value = -448996303.2290542
$ws.Range("D13").Value = 53
$ws.Range("E13").Value = "T"
$ws.Range("B14").Value = 165417181.4159042
$ws.Range("C14").Value = -393879767.4921157
$ws.Range("D14").Value = 53
$ws.Range("E14").Value = "T"
$ws.Range("B15").Value = 179797278.9888872
$ws.Range("C15").Value = -385339572.3715404
$ws.Range("D15").Value = 53
$ws.Range("E15").Value = "T"
$ws.Range("B16").Value = 227210100.4166648
$ws.Range("C16").Value = -464887206.9620859
$ws.Range("D16").Value = 53
$ws.Range("E16").Value = "T"
$ws.Range("B17").Value = 212055122.9125671
$ws.Range("C17").Value = -453803865.4495928
$ws.Range("D17").Value = 53
$ws.Range("E17").Value = "T"
$ws.Range("B18").Value = 219066152.3479342
$ws.Range("C18").Value = -470844080.385905
$ws.Range("D18").Value = 53
$ws.Range("E18").Value = "T"
$ws.Range("B19").Value = 229520900.4912613
$ws.Range("C19").Value = -461115703.439372
$ws.Range("D19").Value = 53
$ws.Range("E19").Value = "T"
$ws.Range("B20").Value = 170189121.3721983
$ws.Range("C20").Value = -402456658.3516716
$ws.Range("D20").Value = 53
$ws.Range("E20").Value = "T"
$ws.Range("B21").Value = 232646418.6580028
$ws.Range("C21").Value = -467100958.0807744
$ws.Range("D21").Value = 53
$ws.Range("E21").Value = "T"
$ws.Range("B22").Value = 221685799.896555
$ws.Range("C22").Value = -450486429.1942855
$ws.Range("D22").Value = 53
$ws.Range("E22").Value = "T"
$ws.Range("B23").Value = 170044434.1339617
$ws.Range("C23").Value = -403681533.2859867
$ws.Range("D23").Value = 53
$ws.Range("E23").Value = "T"
$ws.Range("B24").Value = 170455328.0702242
$ws.Range("C24").Value = -404423592.289158
$ws.Range("D24").Value = 53
$ws.Range("E24").Value = "T"
$ws.Range("B25").Value = 220587369.5919776
$ws.Range("C25").Value = -459065039.3759356
$ws.Range("D25").Value = 53
$ws.Range("E25").Value = "T"
$ws.Range("B26").Value = 204294828.1179895
$ws.Range("C26").Value = -439678649.7171889
$ws.Range("D26").Value = 53
$ws.Range("E26").Value = "T"
$ws.Range("B27").Value = 169897527.0426158
$ws.Range("C27").Value = -402644194.6175289
$ws.Range("D27").Value = 53
$ws.Range("E27").Value = "T"
$ws.Range("B28").Value = 209752092.0299546
$ws.Range("C28").Value = -429279032.9395711
$ws.Range("D28").Value = 53
$ws.Range("E28").Value = "T"
$ws.Range("B29").Value = 181148880.3769667
$ws.Range("C29").Value = -393871673.9209508
$ws.Range("D29").Value = 53
$ws.Range("E29").Value = "T"
$ws.Range("B30").Value = 223592087.1595653
$ws.Range("C30").Value = -454483925.6360297
$ws.Range("D30").Value = 53
$ws.Range("E30").Value = "T"
$ws.Range("B31").Value = 220591388.8574566
$ws.Range("C31").Value = -459070734.1471248
$ws.Range("D31").Value = 53
$ws.Range("E31").Value = "T"
$ws.Range("B32").Value = 169837139.3715255
$ws.Range("C32").Value = -403010049.6333153
$ws.Range("D32").Value = 53
$ws.Range("E32").Value = "T"
$ws.Range("B33").Value = 217075885.120332
$ws.Range("C33").Value = -469975118.9094672
$ws.Range("D33").Value = 53
$ws.Range("E33").Value = "T"
$ws.Range("B34").Value = 170221339.6758617
$ws.Range("C34").Value = -404239109.246271
$ws.Range("D34").Value = 53
$ws.Range("E34").Value = "T"
$ws.Range("B35").Value = 228674745.6979045
$ws.Range("C35").Value = -459433627.7697359
$ws.Range("D35").Value = 53
$ws.Range("E35").Value = "T"
$ws.Range("B36").Value = 226511925.3198773
$ws.Range("C36").Value = -464173405.9727495
$ws.Range("D36").Value = 53
$ws.Range("E36").Value = "T"
$ws.Range("B37").Value = 182849971.5594774
$ws.Range("C37").Value = -397265861.4077961
$ws.Range("D37").Value = 53
$ws.Range("E37").Value = "T"
$ws.Range("B38").Value = 185605712.1449671
$ws.Range("C38").Value = -399507878.409721
$ws.Range("D38").Value = 53
$ws.Range("E38").Value = "T"
$ws.Range("B39").Value = 165158688.1655019
$ws.Range("C39").Value = -397520250.4713072
$ws.Range("D39").Value = 53
$ws.Range("E39").Value = "T"
$ws.Range("B40").Value = 169255388.5440889
$ws.Range("C40").Value = -400598228.2076799
$ws.Range("D40").Value = 53
$ws.Range("E40").Value = "T"
$ws.Range("B41").Value = 220469595.8216356
$ws.Range("C41").Value = -458743532.120712
$ws.Range("D41").Value = 53
$ws.Range("E41").Value = "T"
$ws.Range("B42").Value = 168088018.066688
$ws.Range("C42").Value = -400784453.2981902
$ws.Range("D42").Value = 53
$ws.Range("E42").Value = "T"
$ws.Range("B43").Value = 167778416.8765537
$ws.Range("C43").Value = -396994263.684976
$ws.Range("D43").Value = 53
$ws.Range("E43").Value = "T"
$ws.Range("B44").Value = 189893264.5149418
$ws.Range("C44").Value = -415728052.2039195
$ws.Range("D44").Value = 53
$ws.Range("E44").Value = "T"
$ws.Range("B45").Value = 161294915.6318805
$ws.Range("C45").Value = -386106025.2255547
$ws.Range("D45").Value = 53
$ws.Range("E45").Value = "T"
$ws.Range("B46").Value = 170989661.2193642
$ws.Range("C46").Value = -405566621.7426242
$ws.Range("D46").Value = 53
$ws.Range("E46").Value = "T"
$ws.Range("B47").Value = 167842617.8236445
$ws.Range("C47").Value = -401689953.4288404
$ws.Range("D47").Value = 53
$ws.Range("E47").Value = "T"
$ws.Range("B48").Value = 168348387.5857986
$ws.Range("C48").Value = -398256561.2255399
$ws.Range("D48").Value = 53
$ws.Range("E48").Value = "T"
$ws.Range("B49").Value = 225037801.297646
$ws.Range("C49").Value = -461912078.2331657
$ws.Range("D49").Value = 53
$ws.Range("E49").Value = "T"
$ws.Range("B50").Value = 170869051.6025432
$ws.Range("C50").Value = -403499299.370443
$ws.Range("D50").Value = 53
$ws.Range("E50").Value = "T"
$ws.Range("B51").Value = 193017470.4955127
$ws.Range("C51").Value = -401785101.312815
$ws.Range("D51").Value = 53
$ws.Range("E51").Value = "T"
$ws.Range("B52").Value = 199949949.5776266
$ws.Range("C52").Value = -448945059.420067
$ws.Range("D52").Value = 53
$ws.Range("E52").Value = "T"
$ws.Range("B53").Value = 212560183.8336233
$ws.Range("C53").Value = -460556023.6364397
$ws.Range("D53").Value = 53
$ws.Range("E53").Value = "T"
$ws.Range("B54").Value = 161676484.2880752
$ws.Range("C54").Value = -386240761.038307
$ws.Range("D54").Value = 53
$ws.Range("E54").Value = "T"
$ws.Range("B55").Value = 169788566.1146641
$ws.Range("C55").Value = -401736005.8829324
$ws.Range("D55").Value = 53
$ws.Range("E55").Value = "T"
$ws.Range("B56").Value = 211508905.0259657
$ws.Range("C56").Value = -452247932.9809692
$ws.Range("D56").Value = 53
$ws.Range("E56").Value = "T"
$ws.Range("B57").Value = 159944451.6321533
$ws.Range("C57").Value = -388908680.488093
$ws.Range("D57").Value = 53
$ws.Range("E57").Value = "T"
$ws.Range("B58").Value = 201934275.4551795
$ws.Range("C58").Value = -451332610.6898304
$ws.Range("D58").Value = 53
$ws.Range("E58").Value = "T"
$ws.Range("B59").Value = 184732978.0088556
$ws.Range("C59").Value = -397506742.2204202
$ws.Range("D59").Value = 53
$ws.Range("E59").Value = "T"
$ws.Range("B60").Value = 195705534.4388087
$ws.Range("C60").Value = -440888996.2268443
$ws.Range("D60").Value = 53
$ws.Range("E60").Value = "T"
$ws.Range("B61").Value = 220534136.2407041
$ws.Range("C61").Value = -482430540.6605094
$ws.Range("D61").Value = 53
$ws.Range("E61").Value = "T"
$ws.Range("B62").Value = 225035389.0068487
$ws.Range("C62").Value = -461908760.2337005
$ws.Range("D62").Value = 53
$ws.Range("E62").Value = "T"
$ws.Range("B63").Value = 183320051.9313003
$ws.Range("C63").Value = -390873338.4612513
$ws.Range("D63").Value = 53
$ws.Range("E63").Value = "T"
$ws.Range("B64").Value = 198119290.4729356
$ws.Range("C64").Value = -445812024.3581349
$ws.Range("D64").Value = 53
$ws.Range("E64").Value = "T"
$ws.Range("B65").Value = 183136998.266023
$ws.Range("C65").Value = -425597138.2776437
$ws.Range("D65").Value = 53
$ws.Range("E65").Value = "T"
$ws.Range("B66").Value = 194155942.949601
$ws.Range("C66").Value = -401840817.6504443
$ws.Range("D66").Value = 53
$ws.Range("E66").Value = "T"
$ws.Range("B67").Value = 199951766.5902286
$ws.Range("C67").Value = -448948024.5963792
$ws.Range("D67").Value = 53
$ws.Range("E67").Value = "T"
$ws.Range("B68").Value = 181191745.9251307
$ws.Range("C68").Value = -392441938.8805223
$ws.Range("D68").Value = 53
$ws.Range("E68").Value = "T"
$ws.Range("B69").Value = 181465006.6044621
$ws.Range("C69").Value = -381565732.0586818
$ws.Range("D69").Value = 53
$ws.Range("E69").Value = "T"
$ws.Range("B70").Value = 190126419.3274445
$ws.Range("C70").Value = -405131540.1544722
$ws.Range("D70").Value = 53
$ws.Range("E70").Value = "T"
$ws.Range("B71").Value = 213293560.7974027
$ws.Range("C71").Value = -444046630.9603429
$ws.Range("D71").Value = 53
$ws.Range("E71").Value = "T"
$ws.Range("B72").Value = 170191838.1157665
$ws.Range("C72").Value = -402461543.5417274
$ws.Range("D72").Value = 53
$ws.Range("E72").Value = "T"
$ws.Range("B73").Value = 220588452.1677152
$ws.Range("C73").Value = -459066573.2519769
$ws.Range("D73").Value = 53
$ws.Range("E73").Value = "T"
$ws.Range("B74").Value = 195997045.5293688
$ws.Range("C74").Value = -431996651.364231
$ws.Range("D74").Value = 53
$ws.Range("E74").Value = "T"
$ws.Range("B75").Value = 168450500.9432539
$ws.Range("C75").Value = -401384793.0106061
$ws.Range("D75").Value = 53
$ws.Range("E75").Value = "T"
$ws.Range("B76").Value = 230175350.2550632
$ws.Range("C76").Value = -469965199.3425449
$ws.Range("D76").Value = 53
$ws.Range("E76").Value = "T"
$ws.Range("B77").Value = 208808591.8167622
$ws.Range("C77").Value = -437210044.3358359
$ws.Range("D77").Value = 53
$ws.Range("E77").Value = "T"
$ws.Range("B78").Value = 175243179.3873374
$ws.Range("C78").Value = -381724269.4209624
$ws.Range("D78").Value = 53
$ws.Range("E78").Value = "T"
$ws.Range("B79").Value = 159981689.5436021
$ws.Range("C79").Value = -389034951.2546384
$ws.Range("D79").Value = 53
$ws.Range("E79").Value = "T"
$ws.Range("B80").Value = 201940941.1117887
$ws.Range("C80").Value = -451343395.5166952
$ws.Range("D80").Value = 53
$ws.Range("E80").Value = "T"
$ws.Range("B81").Value = 199038252.9463307
$ws.Range("C81").Value = -447386185.6953286
$ws.Range("D81").Value = 53
$ws.Range("E81").Value = "T"
$ws.Range("B82").Value = 200536386.0293367
$ws.Range("C82").Value = -450110034.1736351
$ws.Range("D82").Value = 53
$ws.Range("E82").Value = "T"
$ws.Range("B83").Value = 205759380.4885022
$ws.Range("C83").Value = -423798822.2369605
$ws.Range("D83").Value = 53
$ws.Range("E83").Value = "T"
$ws.Range("B84").Value = 189088595.2522518
$ws.Range("C84").Value = -395445249.0775359
$ws.Range("D84").Value = 53
$ws.Range("E84").Value = "T"
$ws.Range("B85").Value = 223641556.3890886
$ws.Range("C85").Value = -459278078.5597609
$ws.Range("D85").Value = 53
$ws.Range("E85").Value = "T"
$ws.Range("B86").Value = 201021809.2568428
$ws.Range("C86").Value = -450970544.5014009
$ws.Range("D86").Value = 53
$ws.Range("E86").Value = "T"
$ws.Range("B87").Value = 180714862.6879153
$ws.Range("C87").Value = -386459896.4466831
$ws.Range("D87").Value = 53
$ws.Range("E87").Value = "T"
$ws.Range("B88").Value = 192332759.4216259
$ws.Range("C88").Value = -434773788.9701531
$ws.Range("D88").Value = 53
$ws.Range("E88").Value = "T"
$ws.Range("B89").Value = 160418897.5706784
$ws.Range("C89").Value = -389911002.6655887
$ws.Range("D89").Value = 53
$ws.Range("E89").Value = "T"
$ws.Range("B90").Value = 183326049.5143205
$ws.Range("C90").Value = -425033345.1795313
$ws.Range("D90").Value = 53
$ws.Range("E90").Value = "T"
$ws.Range("B91").Value = 159985281.4310923
$ws.Range("C91").Value = -389041687.233989
$ws.Range("D91").Value = 53
$ws.Range("E91").Value = "T"
$ws.Range("B92").Value = 178908355.1332637
$ws.Range("C92").Value = -404735969.3286093
$ws.Range("D92").Value = 53
$ws.Range("E92").Value = "T"
$ws.Range("B93").Value = 158953626.7530148
$ws.Range("C93").Value = -386767804.1345814
$ws.Range("D93").Value = 53
$ws.Range("E93").Value = "T"
$ws.Range("B94").Value = 159983783.1963968
$ws.Range("C94").Value = -389038877.5547693
$ws.Range("D94").Value = 53
$ws.Range("E94").Value = "T"
$ws.Range("B95").Value = 159031701.4770814
$ws.Range("C95").Value = -387026708.7799101
$ws.Range("D95").Value = 53
$ws.Range("E95").Value = "T"
$ws.Range("B96").Value = 191708627.3984575
$ws.Range("C96").Value = -410724560.5956822
$ws.Range("D96").Value = 53
$ws.Range("E96").Value = "T"
$ws.Range("B97").Value = 164829082.7846155
$ws.Range("C97").Value = -396448779.1607118
$ws.Range("D97").Value = 53
$ws.Range("E97").Value = "T"
$ws.Range("B98").Value = 191830326.5457348
$ws.Range("C98").Value = -420857416.8751051
$ws.Range("D98").Value = 53
$ws.Range("E98").Value = "T"
$ws.Range("B99").Value = 200745874.1465958
$ws.Range("C99").Value = -450730133.2302119
$ws.Range("D99").Value = 53
$ws.Range("E99").Value = "T"
$ws.Range("B100").Value = 195100893.162625
$ws.Range("C100").Value = -440230697.0622374
$ws.Range("D100").Value = 53
$ws.Range("E100").Value = "T"
$ws.Range("B101").Value = 225334535.0004603
$ws.Range("C101").Value = -455721416.1328683
$ws.Range("D101").Value = 53
$ws.Range("E101").Value = "T"
$ws.Range("B102").Value = 170066567.4634582
$ws.Range("C102").Value = -402058122.279879
$ws.Range("D102").Value = 53
$ws.Range("E102").Value = "T"
$ws.Range("B103").Value = 184677066.2989198
$ws.Range("C103").Value = -426343195.7846923
$ws.Range("D103").Value = 53
$ws.Range("E103").Value = "T"
$ws.Range("B104").Value = 190705854.829486
$ws.Range("C104").Value = -429290698.053645
$ws.Range("D104").Value = 53
$ws.Range("E104").Value = "T"
$ws.Range("B105").Value = 173950496.5190559
$ws.Range("C105").Value = -404282480.4928117
$ws.Range("D105").Value = 53
$ws.Range("E105").Value = "T"
$ws.Range("B106").Value = 236488082.5247939
$ws.Range("C106").Value = -502234364.1522337
$ws.Range("D106").Value = 53
$ws.Range("E106").Value = "T"
$ws.Range("B107").Value = 212723687.7217174
$ws.Range("C107").Value = -449317799.2371773
$ws.Range("D107").Value = 53
$ws.Range("E107").Value = "T"
$ws.Range("B108").Value = 169111069.225488
$ws.Range("C108").Value = -402939472.3102995
$ws.Range("D108").Value = 53
$ws.Range("E108").Value = "T"
$ws.Range("B109").Value = 198336644.9666787
$ws.Range("C109").Value = -411866837.4992051
$ws.Range("D109").Value = 53
$ws.Range("E109").Value = "T"
$ws.Range("B110").Value = 158546247.0333468
$ws.Range("C110").Value = -380475198.7775751
$ws.Range("D110").Value = 53
$ws.Range("E110").Value = "T"
$ws.Range("B111").Value = 189729394.5335976
$ws.Range("C111").Value = -403999810.6284639
$ws.Range("D111").Value = 53
$ws.Range("E111").Value = "T"
$ws.Range("B112").Value = 202072665.5388272
$ws.Range("C112").Value = -453520694.1747884
$ws.Range("D112").Value = 53
$ws.Range("E112").Value = "T"
$ws.Range("B113").Value = 161325373.0124752
$ws.Range("C113").Value = -385652928.7783066
$ws.Range("D113").Value = 53
$ws.Range("E113").Value = "T"
$ws.Range("B114").Value = 161054255.9680095
$ws.Range("C114").Value = -385324290.4075406
$ws.Range("D114").Value = 53
$ws.Range("E114").Value = "T"
$ws.Range("B115").Value = 228912456.7268519
$ws.Range("C115").Value = -492802897.0585838
$ws.Range("D115").Value = 53
$ws.Range("E115").Value = "T"
$ws.Range("B116").Value = 189903671.6324413
$ws.Range("C116").Value = -408201000.6104372
$ws.Range("D116").Value = 53
$ws.Range("E116").Value = "T"
$ws.Range("B117").Value = 190128098.1234068
$ws.Range("C117").Value = -415313302.3688121
$ws.Range("D117").Value = 53
$ws.Range("E117").Value = "T"
$ws.Range("B118").Value = 219741114.3635682
$ws.Range("C118").Value = -454396973.9048197
$ws.Range("D118").Value = 53
$ws.Range("E118").Value = "T"
$ws.Range("B119").Value = 191441505.7863299
$ws.Range("C119").Value = -427522610.3461448
$ws.Range("D119").Value = 53
$ws.Range("E119").Value = "T"
$ws.Range("B120").Value = 166744254.6523182
$ws.Range("C120").Value = -395359670.511667
$ws.Range("D120").Value = 53
$ws.Range("E120").Value = "T"
$ws.Range("B121").Value = 159982778.8322647
$ws.Range("C121").Value = -389036994.0382276
$ws.Range("D121").Value = 53
$ws.Range("E121").Value = "T"
$ws.Range("B122").Value = 159450693.4018503
$ws.Range("C122").Value = -382313786.5138384
$ws.Range("D122").Value = 53
$ws.Range("E122").Value = "T"
$ws.Range("B123").Value = 160581153.2411571
$ws.Range("C123").Value = -384341226.8542533
$ws.Range("D123").Value = 53
$ws.Range("E123").Value = "T"
$ws.Range("B124").Value = 204412644.2140465
$ws.Range("C124").Value = -432688320.1416128
$ws.Range("D124").Value = 53
$ws.Range("E124").Value = "T"
$ws.Range("B125").Value = 163617850.8680348
$ws.Range("C125").Value = -388648087.9597438
$ws.Range("D125").Value = 53
$ws.Range("E125").Value = "T"
$ws.Range("B126").Value = 188672409.5696528
$ws.Range("C126").Value = -412151951.1620351
$ws.Range("D126").Value = 53
$ws.Range("E126").Value = "T"
$ws.Range("B127").Value = 180662603.712813
$ws.Range("C127").Value = -380285375.9184469
$ws.Range("D127").Value = 53
$ws.Range("E127").Value = "T"
$ws.Range("B128").Value = 159981312.4340496
$ws.Range("C128").Value = -389034244.04576
$ws.Range("D128").Value = 53
$ws.Range("E128").Value = "T"
$ws.Range("B129").Value = 190553921.2513739
$ws.Range("C129").Value = -418754253.6420302
$ws.Range("D129").Value = 53
$ws.Range("E129").Value = "T"
$ws.Range("B130").Value = 164680010.3926166
$ws.Range("C130").Value = -379685424.6714844
$ws.Range("D130").Value = 53
$ws.Range("E130").Value = "T"
$ws.Range("B131").Value = 180693895.2109168
$ws.Range("C131").Value = -407424798.7042264
$ws.Range("D131").Value = 53
$ws.Range("E131").Value = "T"
$ws.Range("B132").Value = 162283396.6694638
$ws.Range("C132").Value = -386557522.8375871
$ws.Range("D132").Value = 53
$ws.Range("E132").Value = "T"
$ws.Range("B133").Value = 159593752.1135564
$ws.Range("C133").Value = -388307085.0363675
$ws.Range("D133").Value = 53
$ws.Range("E133").Value = "T"
$ws.Range("B134").Value = 183568970.4222474
$ws.Range("C134").Value = -418382782.426956
$ws.Range("D134").Value = 53
$ws.Range("E134").Value = "T"
$ws.Range("B135").Value = 193274830.4161268
$ws.Range("C135").Value = -400423351.3916775
$ws.Range("D135").Value = 53
$ws.Range("E135").Value = "T"
$ws.Range("B136").Value = 169039221.4966475
$ws.Range("C136").Value = -398792434.4795803
$ws.Range("D136").Value = 53
$ws.Range("E136").Value = "T"
$ws.Range("B137").Value = 187767923.9888766
$ws.Range("C137").Value = -412204218.9547121
$ws.Range("D137").Value = 53
$ws.Range("E137").Value = "T"
$ws.Range("B138").Value = 176614885.9983608
$ws.Range("C138").Value = -404281953.2558571
$ws.Range("D138").Value = 53
$ws.Range("E138").Value = "T"
$ws.Range("B139").Value = 214146793.1152762
$ws.Range("C139").Value = -473637174.7999482
$ws.Range("D139").Value = 53
$ws.Range("E139").Value = "T"
$ws.Range("B140").Value = 195694588.3217317
$ws.Range("C140").Value = -437953639.0172651
$ws.Range("D140").Value = 53
$ws.Range("E140").Value = "T"
$ws.Range("B141").Value = 161681400.0862358
$ws.Range("C141").Value = -386249778.728037
$ws.Range("D141").Value = 53
$ws.Range("E141").Value = "T"
$ws.Range("B142").Value = 199808924.882672
$ws.Range("C142").Value = -421352675.5362641
$ws.Range("D142").Value = 53
$ws.Range("E142").Value = "T"
$ws.Range("B143").Value = 187914817.4255285
$ws.Range("C143").Value = -411011705.6538569
$ws.Range("D143").Value = 53
$ws.Range("E143").Value = "T"
$ws.Range("B144").Value = 170030570.7160479
$ws.Range("C144").Value = -373916385.7070824
$ws.Range("D144").Value = 53
$ws.Range("E144").Value = "T"
$ws.Range("B145").Value = 218929775.7707652
$ws.Range("C145").Value = -440791059.8359208
$ws.Range("D145").Value = 53
$ws.Range("E145").Value = "T"
$ws.Range("B146").Value = 162708234.4723862
$ws.Range("C146").Value = -388472907.154898
$ws.Range("D146").Value = 53
$ws.Range("E146").Value = "T"
$ws.Range("B147").Value = 193098927.3384862
$ws.Range("C147").Value = -436513085.7480596
$ws.Range("D147").Value = 53
$ws.Range("E147").Value = "T"
$ws.Range("B148").Value = 230555666.8963441
$ws.Range("C148").Value = -459769917.9560463
$ws.Range("D148").Value = 53
$ws.Range("E148").Value = "T"
$ws.Range("B149").Value = 193727869.7420573
$ws.Range("C149").Value = -433795557.1452767
$ws.Range("D149").Value = 53
$ws.Range("E149").Value = "T"
$ws.Range("B150").Value = 172808170.7906583
$ws.Range("C150").Value = -405162410.6592178
$ws.Range("D150").Value = 53
$ws.Range("E150").Value = "T"
$ws.Range("B151").Value = 190025554.5393861
$ws.Range("C151").Value = -418302784.1040679
$ws.Range("D151").Value = 53
$ws.Range("E151").Value = "T"
$ws.Range("B152").Value = 233758818.4293233
$ws.Range("C152").Value = -497716248.5745634
$ws.Range("D152").Value = 53
$ws.Range("E152").Value = "T"
$ws.Range("B153").Value = 197999048.2211913
$ws.Range("C153").Value = -442524433.6819119
$ws.Range("D153").Value = 53
$ws.Range("E153").Value = "T"
$ws.Range("B154").Value = 170274010.9602531
$ws.Range("C154").Value = -401599067.683601
$ws.Range("D154").Value = 53
$ws.Range("E154").Value = "T"
$ws.Range("B155").Value = 217091423.600224
$ws.Range("C155").Value = -452864053.2432491
$ws.Range("D155").Value = 53
$ws.Range("E155").Value = "T"
$ws.Range("B156").Value = 166872681.8275874
$ws.Range("C156").Value = -394659380.2129357
$ws.Range("D156").Value = 53
$ws.Range("E156").Value = "T"
$ws.Range("B157").Value = 208841063.2337069
$ws.Range("C157").Value = -433398500.873895
$ws.Range("D157").Value = 53
$ws.Range("E157").Value = "T"
$ws.Range("B158").Value = 191697416.7709543
$ws.Range("C158").Value = -400171482.9756765
$ws.Range("D158").Value = 53
$ws.Range("E158").Value = "T"
$ws.Range("B159").Value = 163624090.056206
$ws.Range("C159").Value = -390320148.7865415
$ws.Range("D159").Value = 53
$ws.Range("E159").Value = "T"
$ws.Range("B160").Value = 230432261.7004661
$ws.Range("C160").Value = -459446216.7997958
$ws.Range("D160").Value = 53
$ws.Range("E160").Value = "T"
$ws.Range("B161").Value = 213054465.7445063
$ws.Range("C161").Value = -433863851.7112849
$ws.Range("D161").Value = 53
$ws.Range("E161").Value = "T"
$ws.Range("B162").Value = 192706308.6065789
$ws.Range("C162").Value = -435326852.917184
$ws.Range("D162").Value = 53
$ws.Range("E162").Value = "T"
$ws.Range("B163").Value = 208919250.3967201
$ws.Range("C163").Value = -435841706.9338781
$ws.Range("D163").Value = 53
$ws.Range("E163").Value = "T"
$ws.Range("B164").Value = 169634345.0394309
$ws.Range("C164").Value = -400688439.0056636
$ws.Range("D164").Value = 53
$ws.Range("E164").Value = "T"
$ws.Range("B165").Value = 214076451.9032896
$ws.Range("C165").Value = -434970040.9152226
$ws.Range("D165").Value = 53
$ws.Range("E165").Value = "T"
$ws.Range("B166").Value = 226014241.3440962
$ws.Range("C166").Value = -460462448.9135036
$ws.Range("D166").Value = 53
$ws.Range("E166").Value = "T"
$ws.Range("B167").Value = 169510653.6271445
$ws.Range("C167").Value = -400288688.0457355
$ws.Range("D167").Value = 53
$ws.Range("E167").Value = "T"
$ws.Range("B168").Value = 197544351.2653601
$ws.Range("C168").Value = -428644976.7766817
$ws.Range("D168").Value = 53
$ws.Range("E168").Value = "T"
$ws.Range("B169").Value = 225147600.0272056
$ws.Range("C169").Value = -446621699.3602381
$ws.Range("D169").Value = 53
$ws.Range("E169").Value = "T"
$ws.Range("B170").Value = 204452339.7224762
$ws.Range("C170").Value = -427274707.5131188
$ws.Range("D170").Value = 53
$ws.Range("E170").Value = "T"
$ws.Range("B171").Value = 211436154.0418527
$ws.Range("C171").Value = -432216727.6757636
$ws.Range("D171").Value = 53
$ws.Range("E171").Value = "T"
$ws.Range("B172").Value = 224428130.6488923
$ws.Range("C172").Value = -460835459.0009965
$ws.Range("D172").Value = 53
$ws.Range("E172").Value = "T"
$ws.Range("B173").Value = 230867323.3203529
$ws.Range("C173").Value = -459433945.2831934
$ws.Range("D173").Value = 53
$ws.Range("E173").Value = "T"
$ws.Range("B174").Value = 188099519.5105205
$ws.Range("C174").Value = -414814430.0133743
$ws.Range("D174").Value = 53
$ws.Range("E174").Value = "T"
$ws.Range("B175").Value = 226042122.9440573
$ws.Range("C175").Value = -448979737.9827173
$ws.Range("D175").Value = 53
$ws.Range("E175").Value = "T"
$ws.Range("B176").Value = 227661733.7457307
$ws.Range("C176").Value = -466707593.4775785
$ws.Range("D176").Value = 53
$ws.Range("E176").Value = "T"
$ws.Range("B177").Value = 189877414.6763963
$ws.Range("C177").Value = -414587537.4387551
$ws.Range("D177").Value = 53
$ws.Range("E177").Value = "T"
$ws.Range("B178").Value = 194700531.8312197
$ws.Range("C178").Value = -439028991.2625205
$ws.Range("D178").Value = 53
$ws.Range("E178").Value = "T"
$ws.Range("B179").Value = 193084589.4705866
$ws.Range("C179").Value = -435886815.5279012
$ws.Range("D179").Value = 53
$ws.Range("E179").Value = "T"
$ws.Range("B180").Value = 215805372.3812869
$ws.Range("C180").Value = -465711154.7927014
$ws.Range("D180").Value = 53
$ws.Range("E180").Value = "T"
